$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.995.06"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.909.03"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7908"
$ws.Range("E5").Value = "  +6.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.63"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3162"
$ws.Range("E8").Value = "  +3.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.33"
$ws.Range("E9").Value = "  +3.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06887"
$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08000"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").Value = "1.906.13"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7422"
$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.195"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.06"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "30.000.27"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.92"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.869"
$ws.Range("E18").Value = "  -4.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.63"
$ws.Range("E19").Value = "  +3.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007736"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").Value = "2.143.29"
$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.835"
$ws.Range("E24").Value = "  -3.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.34"
$ws.Range("E25").Value = "  +0.96%  "

$ws.Range("E26").Value = "  -0.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1390"
$ws.Range("E27").Value = "  +10.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.89"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.034"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.367"
$ws.Range("E30").Value = "  +1.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.517"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.313"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.080"
$ws.Range("E33").Value = "  +1.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05507"
$ws.Range("E34").Value = "  +2.31%  "

$ws.Range("E35").Value = "  -2.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7337"
$ws.Range("E36").Value = "  -0.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01925"
$ws.Range("E38").Value = "  -0.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.788"
$ws.Range("E39").Value = "  +0.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.142"
$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.26"
$ws.Range("E42").Value = "  -0.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8372"
$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.872"
$ws.Range("E45").Value = "  -3.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.48"
$ws.Range("E46").Value = "  -1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.547"
$ws.Range("E47").Value = "  -2.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "986.90"
$ws.Range("E48").Value = "  +8.69%  "

$ws.Range("D49").Value = "2.052.81"
$ws.Range("E49").Value = "  -0.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.20"
$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.479"
$ws.Range("E51").Value = "  +0.48%  "
